# The shared "get" routine used to write survey source file paths using the
# old Output\Output Files folder, which caused duplicate files. It now points
# at a "surveys" folder instead, using forward slashes for the directory part
# (retaining a single backslash right before the file name).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPrefix = "C:\Users\paul.jones\Documents\GitHub\ModiffyEfficiency\ModiffyEfficiency\Output\Output Files\"
$newPrefix = "C:/Users/paul.jones/Documents/GitHub/ModiffyEfficiency/ModiffyEfficiency/surveys\"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.StartsWith($oldPrefix)) {
        $cell.Value2 = $newPrefix + $val.Substring($oldPrefix.Length)
    }
}
